# COP25_man_serbia.xlsx — "Final man. annotated files"
#
# The sheet originally had two extra classification columns ("Scale" and
# "Time") between "Unit" and "Principle", and every annotation in column B
# was the literal text "no". The final edit:
#   - drops the "Scale" / "Time" columns entirely
#   - adds a new "Shape" classification column in their place (still column E)
#   - turns every column-B annotation into a numeric 0/1 flag

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Replace the column B text annotations with numeric 0/1 flags ------
$bFlags = @{
    2  = 0
    3  = 0
    4  = 0
    5  = 0
    6  = 0
    7  = 0
    8  = 0
    9  = 0
    10 = 0
    11 = 0
    12 = 1
    13 = 1
    14 = 0
    15 = 0
    16 = 0
    17 = 1
    18 = 1
    19 = 1
}
foreach ($r in $bFlags.Keys) {
    $ws.Cells.Item($r, 2).Value = $bFlags[$r]
}

# --- 2. Drop the "Scale" (E) and "Time" (F) columns -----------------------
$ws.Range("E1:F1").EntireColumn.Delete()

# --- 3. Insert a fresh column in their place and label it "Shape" ---------
$ws.Columns("E").Insert()
$ws.Range("E1").Value = "Shape"

# --- 4. Match the author's final selection (cell B19) ----------------------
[void]$ws.Range("B19").Select()

Write-Output "applied COP25 Serbia annotation edits"
